# Updates market-price derived columns (H-N: currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed market data
# pulled by the scheduled runner. Values are plain numbers (no formulas).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H5").Value = 918
$ws.Range("I5").Value = 339.75
$ws.Range("J5").Value = 1689
$ws.Range("K5").Value = 339.75
$ws.Range("L5").Value = 1689
$ws.Range("M5").Value = -224.75
$ws.Range("N5").Value = -1919

$ws.Range("H18").Value = 1130.75
$ws.Range("I18").Value = 449.42856
$ws.Range("K18").Value = 449.42856
$ws.Range("M18").Value = -165.42856

$ws.Range("H80").Value = 2592.7144
$ws.Range("I80").Value = 1159.6
$ws.Range("J80").Value = 3040.5625
$ws.Range("K80").Value = 3478.8
$ws.Range("L80").Value = 9121.6875
$ws.Range("M80").Value = -2480.8
$ws.Range("N80").Value = -11117.6875

$ws.Range("H83").Value = 2592.7144
$ws.Range("I83").Value = 1159.6
$ws.Range("J83").Value = 3040.5625
$ws.Range("K83").Value = 10436.4
$ws.Range("L83").Value = 27365.0625
$ws.Range("M83").Value = -5444.4
$ws.Range("N83").Value = -37349.0625

$ws.Range("H88").Value = 4002.6
$ws.Range("J88").Value = 4002.6
$ws.Range("L88").Value = 4002.6
$ws.Range("N88").Value = -4814.6

$ws.Range("H91").Value = 4002.6
$ws.Range("J91").Value = 4002.6
$ws.Range("L91").Value = 4002.6
$ws.Range("N91").Value = -6810.6

$ws.Range("H129").Value = 14894.4375
$ws.Range("I129").Value = 1731.5
$ws.Range("J129").Value = 22792.2
$ws.Range("K129").Value = 5194.5
$ws.Range("L129").Value = 68376.60000000001
$ws.Range("M129").Value = -194.5
$ws.Range("N129").Value = -78376.60000000001

$ws.Range("H138").Value = 5063.4614
$ws.Range("J138").Value = 7920
$ws.Range("L138").Value = 23760
$ws.Range("N138").Value = -34040

$ws.Range("H141").Value = 16405.7
$ws.Range("I141").Value = 25974.75
$ws.Range("J141").Value = 10026.333
$ws.Range("K141").Value = 77924.25
$ws.Range("L141").Value = 30078.999
$ws.Range("M141").Value = -72744.25
$ws.Range("N141").Value = -40438.999


$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H5").Value = 423.91666
$ws.Range("I5").Value = 439.72726
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 439.72726
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = -327.72726
$ws.Range("N5").Value = -474

$ws.Range("H32").Value = 3674.4658
$ws.Range("I32").Value = 3102.0896
$ws.Range("K32").Value = 3102.0896
$ws.Range("M32").Value = -2815.0896

$ws.Range("H63").Value = 2000
$ws.Range("J63").Value = 2000
$ws.Range("L63").Value = 2000
$ws.Range("N63").Value = -3372

$ws.Range("H66").Value = 2000
$ws.Range("J66").Value = 2000
$ws.Range("L66").Value = 10000
$ws.Range("N66").Value = -16864


$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H4").Value = 423.91666
$ws.Range("I4").Value = 439.72726
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 439.72726
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = -324.72726
$ws.Range("N4").Value = -480

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H86").Value = 5349.2046
$ws.Range("I86").Value = 6147.2964
$ws.Range("J86").Value = 4081.647
$ws.Range("K86").Value = 6147.2964
$ws.Range("L86").Value = 4081.647
$ws.Range("M86").Value = -5024.2964
$ws.Range("N86").Value = -6327.647

$ws.Range("H89").Value = 5349.2046
$ws.Range("I89").Value = 6147.2964
$ws.Range("J89").Value = 4081.647
$ws.Range("K89").Value = 30736.482
$ws.Range("L89").Value = 20408.235
$ws.Range("M89").Value = -25120.482
$ws.Range("N89").Value = -31640.235

$ws.Range("H94").Value = 1746.6471
$ws.Range("I94").Value = 1510.7
$ws.Range("J94").Value = 2083.7144
$ws.Range("K94").Value = 1510.7
$ws.Range("L94").Value = 2083.7144
$ws.Range("M94").Value = -1059.7
$ws.Range("N94").Value = -2985.7144

$ws.Range("H105").Value = 8697.956
$ws.Range("I105").Value = 24964
$ws.Range("K105").Value = 24964
$ws.Range("M105").Value = -23217

$ws.Range("H134").Value = 3941.0293
$ws.Range("I134").Value = 4149.3706
$ws.Range("J134").Value = 3137.4285
$ws.Range("K134").Value = 12448.1118
$ws.Range("L134").Value = 9412.2855
$ws.Range("M134").Value = -9913.111800000001
$ws.Range("N134").Value = -14482.2855


$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H7").Value = 249.20833
$ws.Range("I7").Value = 131.77777
$ws.Range("J7").Value = 319.66666
$ws.Range("K7").Value = 131.77777
$ws.Range("L7").Value = 319.66666
$ws.Range("M7").Value = -18.77777
$ws.Range("N7").Value = -545.66666

$ws.Range("H122").Value = 2106893.5
$ws.Range("I122").Value = 1412.6842
$ws.Range("J122").Value = 8774250
$ws.Range("K122").Value = 4238.0526
$ws.Range("L122").Value = 26322750
$ws.Range("M122").Value = -1788.0526
$ws.Range("N122").Value = -26327650


$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H107").Value = 1587.2
$ws.Range("J107").Value = 1704.3334
$ws.Range("L107").Value = 5113.0002
$ws.Range("N107").Value = -8953.0002

$ws.Range("H132").Value = 3926144.5
$ws.Range("I132").Value = 1786.625
$ws.Range("K132").Value = 16079.625
$ws.Range("M132").Value = -13549.625


$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H132").Value = 5044.2856
$ws.Range("I132").Value = 4961.875
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 14885.625
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -12355.625
$ws.Range("N132").Value = -32060


$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 5244.75
$ws.Range("I40").Value = 5294.1665
$ws.Range("K40").Value = 5294.1665
$ws.Range("M40").Value = -5158.1665

$ws.Range("H46").Value = 1472.1852
$ws.Range("I46").Value = 1032.5652
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 1032.5652
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -844.5652
$ws.Range("N46").Value = -4376

$ws.Range("H93").Value = 661256.4
$ws.Range("J93").Value = 1896168.8
$ws.Range("L93").Value = 1896168.8
$ws.Range("N93").Value = -1898664.8


$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H6").Value = 5000
$ws.Range("J6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5230

